$wb = $excel.ActiveWorkbook

# Sheet ALC, row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1253310.6
$ws.Range("J19").Value = 4697.8
$ws.Range("L19").Value = 4697.8
$ws.Range("N19").Value = -5047.8

# Sheet ALC, row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 125
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 125
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 125
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -583

# Sheet ALC, row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 15616.286
$ws.Range("I80").Value = 2117.6667
$ws.Range("J80").Value = 25740.25
$ws.Range("K80").Value = 6353.000100000001
$ws.Range("L80").Value = 77220.75
$ws.Range("M80").Value = -5355.000100000001
$ws.Range("N80").Value = -79216.75

# Sheet ALC, row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 15616.286
$ws.Range("I83").Value = 2117.6667
$ws.Range("J83").Value = 25740.25
$ws.Range("K83").Value = 19059.0003
$ws.Range("L83").Value = 231662.25
$ws.Range("M83").Value = -14067.0003
$ws.Range("N83").Value = -241646.25

# Sheet ALC, row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 879739.5
$ws.Range("I92").Value = 1119486
$ws.Range("K92").Value = 1119486
$ws.Range("M92").Value = -1118238

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1214.9615
$ws.Range("J132").Value = 899.5
$ws.Range("L132").Value = 2698.5
$ws.Range("N132").Value = -7758.5

# Sheet ALC, row 134
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 44607.69
$ws.Range("J134").Value = 44607.69
$ws.Range("L134").Value = 44607.69
$ws.Range("N134").Value = -54747.69

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2516.0977
$ws.Range("I138").Value = 2398.5217
$ws.Range("K138").Value = 7195.5651
$ws.Range("M138").Value = -2055.5651

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2584755.2
$ws.Range("I2").Value = 2907725
$ws.Range("K2").Value = 2907725
$ws.Range("M2").Value = -2907612

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10347.272
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 10347.272
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 10347.272
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -10921.272

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4457.15
$ws.Range("I61").Value = 2778.3125
$ws.Range("J61").Value = 11172.5
$ws.Range("K61").Value = 2778.3125
$ws.Range("L61").Value = 11172.5
$ws.Range("M61").Value = -2566.3125
$ws.Range("N61").Value = -11596.5

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 811.3
$ws.Range("I74").Value = 794.4483
$ws.Range("K74").Value = 794.4483
$ws.Range("M74").Value = 79.55169999999998

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 811.3
$ws.Range("I77").Value = 794.4483
$ws.Range("K77").Value = 3972.2415
$ws.Range("M77").Value = 395.7584999999999

# Sheet ARM, row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 619.7
$ws.Range("I97").Value = 577.44446
$ws.Range("K97").Value = 577.44446
$ws.Range("M97").Value = -81.44446000000005

# Sheet ARM, row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2035.381
$ws.Range("I102").Value = 1813.5
$ws.Range("J102").Value = 3366.6667
$ws.Range("K102").Value = 1813.5
$ws.Range("L102").Value = 3366.6667
$ws.Range("M102").Value = -191.5
$ws.Range("N102").Value = -6610.6667

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2584755.2
$ws.Range("I116").Value = 2907725
$ws.Range("K116").Value = 2907725
$ws.Range("M116").Value = -2905431

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4457.15
$ws.Range("I136").Value = 2778.3125
$ws.Range("J136").Value = 11172.5
$ws.Range("K136").Value = 8334.9375
$ws.Range("L136").Value = 33517.5
$ws.Range("M136").Value = -5784.9375
$ws.Range("N136").Value = -38617.5

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2584755.2
$ws.Range("I3").Value = 2907725
$ws.Range("K3").Value = 2907725
$ws.Range("M3").Value = -2907611

# Sheet BSM, row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1769.2
$ws.Range("I20").Value = 1995.1818
$ws.Range("J20").Value = 1147.75
$ws.Range("K20").Value = 1995.1818
$ws.Range("L20").Value = 1147.75
$ws.Range("M20").Value = -1748.1818
$ws.Range("N20").Value = -1641.75

# Sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1847.3846
$ws.Range("I107").Value = 1092.3636
$ws.Range("K107").Value = 1092.3636
$ws.Range("M107").Value = 827.6364000000001

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7077.357
$ws.Range("I134").Value = 6681.8057
$ws.Range("J134").Value = 9450.666999999999
$ws.Range("K134").Value = 20045.4171
$ws.Range("L134").Value = 28352.001
$ws.Range("M134").Value = -17510.4171
$ws.Range("N134").Value = -33422.001

# Sheet CRP, row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 665
$ws.Range("I16").Value = 610
$ws.Range("J16").Value = 995
$ws.Range("K16").Value = 610
$ws.Range("L16").Value = 995
$ws.Range("M16").Value = -323
$ws.Range("N16").Value = -1569

# Sheet CRP, row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1366.6666
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 1580
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 1580
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = -2280

# Sheet CRP, row 23
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 65766
$ws.Range("I23").Value = 49800
$ws.Range("K23").Value = 49800
$ws.Range("M23").Value = -49560

# Sheet CRP, row 27
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H27").Value = 65766
$ws.Range("I27").Value = 49800
$ws.Range("K27").Value = 49800
$ws.Range("M27").Value = -49608

# Sheet CRP, row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 30000
$ws.Range("J74").Value = 30000
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -31748

# Sheet CRP, row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 30000
$ws.Range("J77").Value = 30000
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -98736

# Sheet CRP, row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 9999999
$ws.Range("I99").Value = 9999999
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 9999999
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -9998501
$ws.Range("N99").ClearContents()

# Sheet CRP, row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 665
$ws.Range("I113").Value = 610
$ws.Range("J113").Value = 995
$ws.Range("K113").Value = 610
$ws.Range("L113").Value = 995
$ws.Range("M113").Value = 1560
$ws.Range("N113").Value = -5335

# Sheet CRP, row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 9999999
$ws.Range("I126").Value = 9999999
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 29999997
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -29997527
$ws.Range("N126").ClearContents()

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2875.4707
$ws.Range("I132").Value = 1289.5
$ws.Range("K132").Value = 3868.5
$ws.Range("M132").Value = -1338.5

# Sheet CUL, row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 6504.091
$ws.Range("I55").Value = 50252
$ws.Range("J55").Value = 2129.3
$ws.Range("K55").Value = 150756
$ws.Range("L55").Value = 6387.900000000001
$ws.Range("M55").Value = -150579
$ws.Range("N55").Value = -6741.900000000001

# Sheet CUL, row 88
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 6333.3335
$ws.Range("I88").Value = 3000
$ws.Range("J88").Value = 8000
$ws.Range("K88").Value = 9000
$ws.Range("L88").Value = 24000
$ws.Range("M88").Value = -8572
$ws.Range("N88").Value = -24856

# Sheet CUL, row 91
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 6333.3335
$ws.Range("I91").Value = 3000
$ws.Range("J91").Value = 8000
$ws.Range("K91").Value = 9000
$ws.Range("L91").Value = 24000
$ws.Range("M91").Value = -7518
$ws.Range("N91").Value = -26964

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 815.2857
$ws.Range("I131").Value = 531.8333
$ws.Range("J131").Value = 833.7717
$ws.Range("K131").Value = 1595.4999
$ws.Range("L131").Value = 2501.3151
$ws.Range("M131").Value = 3444.5001
$ws.Range("N131").Value = -12581.3151

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2587.7222
$ws.Range("I102").Value = 2726.6
$ws.Range("K102").Value = 2726.6
$ws.Range("M102").Value = -1104.6

# Sheet LTW, row 106
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 17256.25
$ws.Range("J106").Value = 17256.25
$ws.Range("L106").Value = 17256.25
$ws.Range("N106").Value = -19780.25

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1722.7059
$ws.Range("I132").Value = 1461.7894
$ws.Range("K132").Value = 4385.3682
$ws.Range("M132").Value = -1855.3682

# Sheet WVR, row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 4000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -6122

# Sheet WVR, row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 20000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -30608

# Sheet WVR, row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 7367.724
$ws.Range("J126").Value = 8143.6665
$ws.Range("L126").Value = 24430.9995
$ws.Range("N126").Value = -29370.9995

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2366
$ws.Range("I132").Value = 2055.6924
$ws.Range("K132").Value = 6167.0772
$ws.Range("M132").Value = -3637.0772
